# Convention change to support multi-axle vehicles.
# Adds two new vehicle sheets ("Truck_Amandla" and "Trailer_Kumanzi"),
# both cloned from the existing "Trailer_Thwala" sheet, with the
# Truck sheet positioned before Trailer_Thwala and the new trailer
# sheet positioned after it.

$wb = $excel.ActiveWorkbook

# --- Build "Truck_Amandla" as a copy of Trailer_Thwala, placed BEFORE it ---
$thwala = $wb.Worksheets.Item("Trailer_Thwala")
$thwala.Copy($thwala, $null)
$truck = $wb.ActiveSheet
$truck.Name = "Truck_Amandla"

# Trim the two trailing blank rows (10 -> 8 rows of data)
$truck.Rows("9:10").Delete()

# Update the CAD instance labels
$truck.Range("H3").Value = "CAD_Truck_Amandla"
$truck.Range("H4").Value = "CAD_Truck_Amandla"

# Update the sOffset (row 7) and Opacity (row 8) values
$truck.Range("F7").Value = 0.6
$truck.Range("G7").Value = 0.8
$truck.Range("H7").Value = 1
$truck.Range("H8").Value = 1

$truck.Range("G23").Select() | Out-Null

# --- Build "Trailer_Kumanzi" as a copy of Trailer_Thwala, placed AFTER it ---
$thwala2 = $wb.Worksheets.Item("Trailer_Thwala")
$thwala2.Copy($null, $thwala2)
$kumanzi = $wb.ActiveSheet
$kumanzi.Name = "Trailer_Kumanzi"

# Trim the two trailing blank rows (10 -> 8 rows of data)
$kumanzi.Rows("9:10").Delete()

# Update the CAD instance labels
$kumanzi.Range("H3").Value = "CAD_Trailer_Kumanzi"
$kumanzi.Range("H4").Value = "CAD_Trailer_Kumanzi"

# Update the sOffset (row 7) and Opacity (row 8) values
$kumanzi.Range("F7").Value = 1
$kumanzi.Range("G7").Value = 0.75
$kumanzi.Range("H7").Value = 0.055
$kumanzi.Range("H8").Value = 0.5

$kumanzi.Range("H8").Select() | Out-Null

# Trailer_Kumanzi is the sheet that ends up active/selected
$kumanzi.Activate()
